$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; this shifts existing rows 34-152 down to 35-153
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the new data record
$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44600
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = 100112029
$ws.Range("G34").Value = "Orégano"
$ws.Range("H34").Value = "Sin especificar"
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 30
$ws.Range("K34").Value = 8000
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = 8667
$ws.Range("N34").Value = "`$/docena de atados"
$ws.Range("O34").Value = "Región Metropolitana"
$ws.Range("P34").Value = 2889
$ws.Range("Q34").Value = 3
$ws.Range("R34").Value = "Hortaliza"
